$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 50011000
$ws.Range("I62").Value = 66678000
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 66678000
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -66677376
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 50011000
$ws.Range("I65").Value = 66678000
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 333390000
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -333386880
$ws.Range("N65").Value = -56240
$ws.Range("H107").Value = 896.1579
$ws.Range("I107").Value = 884.6429000000001
$ws.Range("K107").Value = 884.6429000000001
$ws.Range("M107").Value = 1035.3571
$ws.Range("H132").Value = 9010801
$ws.Range("I132").Value = 10754060
$ws.Range("K132").Value = 32262180
$ws.Range("M132").Value = -32259650
$ws.Range("H133").Value = 74279.14
$ws.Range("J133").Value = 74279.14
$ws.Range("L133").Value = 74279.14
$ws.Range("N133").Value = -84399.14
$ws.Range("H136").Value = 45555
$ws.Range("J136").Value = 45555
$ws.Range("L136").Value = 45555
$ws.Range("N136").Value = -55755
$ws.Range("H137").Value = 7616.0586
$ws.Range("J137").Value = 34799.668
$ws.Range("L137").Value = 104399.004
$ws.Range("N137").Value = -109499.004
$ws.Range("H138").Value = 20510.209
$ws.Range("J138").Value = 5337.317
$ws.Range("L138").Value = 16011.951
$ws.Range("N138").Value = -26291.951

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3749.6047
$ws.Range("I32").Value = 3447.3901
$ws.Range("K32").Value = 3447.3901
$ws.Range("M32").Value = -3160.3901
$ws.Range("H97").Value = 3150.25
$ws.Range("I97").Value = 2430.3
$ws.Range("K97").Value = 2430.3
$ws.Range("M97").Value = -1934.3
$ws.Range("H110").Value = 7480.6177
$ws.Range("I110").Value = 10054.6
$ws.Range("J110").Value = 3803.5
$ws.Range("K110").Value = 10054.6
$ws.Range("L110").Value = 3803.5
$ws.Range("M110").Value = -8009.6
$ws.Range("N110").Value = -7893.5
$ws.Range("H132").Value = 1735.171
$ws.Range("I132").Value = 1472.5
$ws.Range("K132").Value = 4417.5
$ws.Range("M132").Value = -1887.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3556.8
$ws.Range("I20").Value = 2646.3845
$ws.Range("J20").Value = 4543.0835
$ws.Range("K20").Value = 2646.3845
$ws.Range("L20").Value = 4543.0835
$ws.Range("M20").Value = -2399.3845
$ws.Range("N20").Value = -5037.0835
$ws.Range("H56").Value = 31650
$ws.Range("J56").Value = 32475
$ws.Range("L56").Value = 32475
$ws.Range("N56").Value = -33953
$ws.Range("H107").Value = 1970.4062
$ws.Range("J107").Value = 3600
$ws.Range("L107").Value = 3600
$ws.Range("N107").Value = -7440
$ws.Range("H134").Value = 1676.6
$ws.Range("I134").Value = 1704.26
$ws.Range("K134").Value = 5112.78
$ws.Range("M134").Value = -2577.78

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23327.438
$ws.Range("I31").Value = 32362.848
$ws.Range("J31").Value = 3449.5334
$ws.Range("K31").Value = 32362.848
$ws.Range("L31").Value = 3449.5334
$ws.Range("M31").Value = -32067.848
$ws.Range("N31").Value = -4039.5334
$ws.Range("H34").Value = 23327.438
$ws.Range("I34").Value = 32362.848
$ws.Range("J34").Value = 3449.5334
$ws.Range("K34").Value = 32362.848
$ws.Range("L34").Value = 3449.5334
$ws.Range("M34").Value = -32160.848
$ws.Range("N34").Value = -3853.5334
$ws.Range("H58").Value = 2935.8333
$ws.Range("I58").Value = 2683.9487
$ws.Range("J58").Value = 3590.7334
$ws.Range("K58").Value = 2683.9487
$ws.Range("L58").Value = 3590.7334
$ws.Range("M58").Value = -2480.9487
$ws.Range("N58").Value = -3996.7334
$ws.Range("H107").Value = 1471.3889
$ws.Range("I107").Value = 2033
$ws.Range("K107").Value = 2033
$ws.Range("M107").Value = -113
$ws.Range("H132").Value = 108610.336
$ws.Range("I132").Value = 130871.125
$ws.Range("K132").Value = 392613.375
$ws.Range("M132").Value = -390083.375
$ws.Range("H134").Value = 20941.56
$ws.Range("I134").Value = 14836.312
$ws.Range("J134").Value = 75888.8
$ws.Range("K134").Value = 44508.936
$ws.Range("L134").Value = 227666.4
$ws.Range("M134").Value = -41973.936
$ws.Range("N134").Value = -232736.4
$ws.Range("H136").Value = 2935.8333
$ws.Range("I136").Value = 2683.9487
$ws.Range("J136").Value = 3590.7334
$ws.Range("K136").Value = 8051.8461
$ws.Range("L136").Value = 10772.2002
$ws.Range("M136").Value = -5501.8461
$ws.Range("N136").Value = -15872.2002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1136.7222
$ws.Range("I132").Value = 1094.2903
$ws.Range("J132").Value = 1399.8
$ws.Range("K132").Value = 9848.6127
$ws.Range("L132").Value = 12598.2
$ws.Range("M132").Value = -7318.6127
$ws.Range("N132").Value = -17658.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 5663.684
$ws.Range("I2").Value = 10321.7
$ws.Range("J2").Value = 488.1111
$ws.Range("K2").Value = 10321.7
$ws.Range("L2").Value = 488.1111
$ws.Range("M2").Value = -10208.7
$ws.Range("N2").Value = -714.1111000000001
$ws.Range("H36").Value = 15083.333
$ws.Range("J36").Value = 15200
$ws.Range("L36").Value = 15200
$ws.Range("N36").Value = -16170
$ws.Range("H40").Value = 34835
$ws.Range("J40").Value = 34835
$ws.Range("L40").Value = 34835
$ws.Range("N40").Value = -35137
$ws.Range("H70").Value = 20778.875
$ws.Range("I70").Value = 5806.5
$ws.Range("J70").Value = 25769.666
$ws.Range("K70").Value = 5806.5
$ws.Range("L70").Value = 25769.666
$ws.Range("M70").Value = -5536.5
$ws.Range("N70").Value = -26309.666
$ws.Range("H73").Value = 20778.875
$ws.Range("I73").Value = 5806.5
$ws.Range("J73").Value = 25769.666
$ws.Range("K73").Value = 5806.5
$ws.Range("L73").Value = 25769.666
$ws.Range("M73").Value = -4870.5
$ws.Range("N73").Value = -27641.666
$ws.Range("H117").Value = 84497.5
$ws.Range("J117").Value = 84497.5
$ws.Range("L117").Value = 84497.5
$ws.Range("N117").Value = -91381.5
$ws.Range("H123").Value = 52933.332
$ws.Range("J123").Value = 52933.332
$ws.Range("L123").Value = 52933.332
$ws.Range("N123").Value = -57833.332
$ws.Range("H132").Value = 2049.7778
$ws.Range("I132").Value = 1883.25
$ws.Range("K132").Value = 5649.75
$ws.Range("M132").Value = -3119.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 9651499
$ws.Range("I23").Value = 6636427
$ws.Range("K23").Value = 6636427
$ws.Range("M23").Value = -6636197
$ws.Range("H33").Value = 18000
$ws.Range("I33").Value = 18000
$ws.Range("K33").Value = 18000
$ws.Range("M33").Value = -17710
$ws.Range("H106").Value = 15249.75
$ws.Range("J106").Value = 15249.75
$ws.Range("L106").Value = 15249.75
$ws.Range("N106").Value = -17773.75
$ws.Range("H122").Value = 6708.2666
$ws.Range("I122").Value = 6804.5713
$ws.Range("K122").Value = 20413.7139
$ws.Range("M122").Value = -17963.7139
$ws.Range("H132").Value = 2407.4036
$ws.Range("I132").Value = 2224.54
$ws.Range("K132").Value = 6673.62
$ws.Range("M132").Value = -4143.62
$ws.Range("H133").Value = 47940
$ws.Range("J133").Value = 60325
$ws.Range("L133").Value = 60325
$ws.Range("N133").Value = -65385
$ws.Range("H136").Value = 48745.953
$ws.Range("I136").Value = 62053.94
$ws.Range("J136").Value = 3498.8
$ws.Range("K136").Value = 186161.82
$ws.Range("L136").Value = 10496.4
$ws.Range("M136").Value = -183611.82
$ws.Range("N136").Value = -15596.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1728.1936
$ws.Range("I122").Value = 1658.5883
$ws.Range("K122").Value = 4975.7649
$ws.Range("M122").Value = -2525.7649
$ws.Range("H132").Value = 1587.6531
$ws.Range("I132").Value = 1448.3414
$ws.Range("K132").Value = 4345.0242
$ws.Range("M132").Value = -1815.0242
$ws.Range("H136").Value = 2851.6316
$ws.Range("I136").Value = 2981.75
$ws.Range("J136").Value = 2628.5715
$ws.Range("K136").Value = 8945.25
$ws.Range("L136").Value = 7885.7145
$ws.Range("M136").Value = -6395.25
$ws.Range("N136").Value = -12985.7145
